# Auto-generated COM-interop script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.462.08"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.807.98"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'225.31"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "'0.587"
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'38.10"
$ws.Range("E8").Value = "  +5.73%  "
$ws.Range("E9").Value = "  -4.24%  "
$ws.Range("D10").Value = "'0.0671"
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").Value = "'0.0973"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "2.070.39"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "'11.09"
$ws.Range("E13").Value = "  -4.95%  "
$ws.Range("D14").Value = "1.811.26"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.454.71"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.628"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").Value = "'4.40"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").Value = "'67.96"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").Value = "'241.93"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "0.0₃0769"
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").Value = "'11.11"
$ws.Range("E21").Value = "  -4.19%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'4.09"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").Value = "'2.19"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").Value = "'169.83"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").Value = "'7.69"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("D27").Value = "'17.53"
$ws.Range("E27").Value = "  +4.16%  "
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "'3.76"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("D32").Value = "'0.0514"
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("D33").Value = "'3.83"
$ws.Range("E33").Value = "  -4.60%  "
$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "1.352.14"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.06"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.639"
$ws.Range("E37").Value = "  -4.69%  "
$ws.Range("D38").Value = "'0.0188"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").Value = "'2.32"
$ws.Range("E39").Value = "  -5.76%  "
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("D41").Value = "'1.22"
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'81.52"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.79"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'0.942"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").Value = "'13.63"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("D47").Value = "1.971.42"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "'5.72"
$ws.Range("E48").Value = "  -5.01%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'102.00"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("E51").Value = "  -4.14%  "
